# "testing with real data": rename header labels on Sheet1 row 1.
#   C1: frequency          -> angular frequency
#   A1: z_real (lowercase)  -> Z_real (capitalized)
#   B1: z_imag (lowercase)  -> Z_imag (capitalized)
# Columns D:G (eff_cap, applied voltage, J_ph, J) and all data rows are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "angular frequency"
$ws.Range("A1").Value = "Z_real"
$ws.Range("B1").Value = "Z_imag"

# Move the sheet's active selection from I5 to B1
$ws.Range("B1").Select()
